# Quarterly indexing esoteric bug-fix operation
# Shift the evaluation rows down by one (row N <- row (N-1)'s data) for rows 3..11,
# and populate row 2 with the newly (re-indexed) computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G")

# Work from the bottom row upward so we never clobber a source row before
# it has been read.
for ($r = 11; $r -ge 3; $r--) {
    $src = $r - 1
    foreach ($col in $cols) {
        $val = $ws.Range("$col$src").Value()
        $ws.Range("$col$r").Value = $val
    }
}

# Row 2 gets the freshly computed values
$ws.Range("B2").Value = -0.02907897629796788
$ws.Range("C2").Value = 0.3131278957257717
$ws.Range("D2").Value = 0.181524606355785
$ws.Range("E2").Value = 0.4260570458938391
$ws.Range("F2").Value = 0.43998257208981
$ws.Range("G2").Value = 15
